$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-tracking entry for 2021-08-03 (row 6): arrival 17:15, departure 20:24.
# The "Temps" (E6) and "Temps total" (K5) formulas already on the sheet recalc automatically.
$ws.Range("B6").Value = 44411
$ws.Range("C6").Value = 0.71875
$ws.Range("D6").Value = 0.85

# Cursor moves down to the next empty row after entering the row.
$ws.Range("B7").Select()
